$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("C2").Value = "The applicant has a score of 72.69 as he possesses skills such as ReactJS, JavaScript, and a strong focus on frontend development, which align well with the job description. However, the lack of experience with NextJS and CSS may have led to a slightly lower score."
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = "The applicant demonstrates good communication, teamwork, problem-solving skills, and a willingness to adapt to a new culture in Japan. The candidate's eagerness to contribute to cutting-edge AI/ML technologies and learn new skills aligns with the company's values, earning a score of 4 out of 5."

# Row 3
$ws.Range("C3").Value = "The applicant has a good score due to possessing skills like CSS, NodeJS, HTML, Socket.IO, WebRTC, ExpressJS, and ReactJS, which are relevant to frontend development. However, the lack of experience with NextJS and JavaScript, as required by the company, might have slightly impacted the score."
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = "The applicant has shown adaptability, teamwork, and willingness to learn and contribute, making him suitable for the company's work culture and potentially working in Japan."

# Row 4
$ws.Range("C4").Value = "The applicant has a high score due to their strong projects involving computer vision, machine learning, and web development, aligning well with the required skills for the Frontend Engineer Intern position. Although lacking direct experience with ReactJS, the applicant's proficiency in related technologies and their comprehensive web development project demonstrate their potential to succeed."
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = "The applicant demonstrates strong adaptability, teamwork, and problem-solving skills, as evidenced by their proactive approach to addressing weaknesses and overcoming challenges. Their willingness and interest in working in Japan, along with their desire to contribute to cutting-edge technologies, reflects a positive attitude and alignment with the company's goals."

$wb.Save()
